$wb = $excel.ActiveWorkbook

# --- Sheet "Summary": update aggregate metrics ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.95   # Current Capital
$summary.Range("B4").Value = -0.05     # Total P&L $
$summary.Range("B5").Value = -0.06     # Total P&L %
$summary.Range("B6").Value = 17        # Total Trades
$summary.Range("B8").Value = 10        # Losing Trades
$summary.Range("B9").Value = 35.29     # Win Rate %

# --- Sheet "Strategy Status": update MarketMaking row (row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.95      # Capital
$status.Range("D4").Value = 17         # Trades
$status.Range("E4").Value = -0.05      # P&L $
$status.Range("F4").Value = -0.05      # P&L %
$status.Range("G4").Value = 35.29      # Win Rate %

# --- Append new trade (Trade #17) to "All Trades" and "MarketMaking" sheets ---
function Add-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value = 17
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 3).Value = "07:54:23"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.02
    $ws.Cells.Item($row, 7).Value = 0.01
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -50
    $ws.Cells.Item($row, 10).Value = -0.01
    $ws.Cells.Item($row, 11).Value = 99.95
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.16
    # Clear the temporary text-number-format so the cell keeps default styling
    # (matches the rest of the sheet, which has no explicit cell style).
    $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 17)).ClearFormats()
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 18

$mm = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $mm 18
